$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-style row 18 to the "section separator" look (copy format from row 15,
#     which already uses styles 8/9 for this pattern) ---
$srcSep = $ws.Range("A15:E15")
$dstSep = $ws.Range("A18:E18")
$srcSep.Copy()
$dstSep.PasteSpecial(-4122)

# --- Add the three new data rows (19, 20, 21) ---
# Values are assigned in the exact order the original author entered them
# (column-wise paste order), since that order determines how new entries
# are appended to the shared-strings table.

$ws.Cells.Item(19, 3).Value = " I thank you sincerely."
$ws.Cells.Item(20, 3).Value = " It was your doing that brought\npeace to our world."
$ws.Cells.Item(21, 3).Value = " The fact that I can lounge on the\nrocks here... That\'s your doing too. Ho-ho-ho!"

$ws.Cells.Item(19, 1).Value = " SCRIPT/P02P01A/us0101.ssb"

$ws.Cells.Item(19, 4).Value = " Я благодарю вас от всего\nсердца."
$ws.Cells.Item(20, 4).Value = " Именно вы смогли даровать покой\nнашему миру."
$ws.Cells.Item(21, 4).Value = " То, что я могу тут лежать на\nкамнях... Это тоже ваша заслуга. Хо-хо-хо!"

$ws.Cells.Item(19, 5).Value = " Ÿ áìàãïäàñý âàò ïó âòåãï\nòåñäøà."
$ws.Cells.Item(20, 5).Value = " Éíåîîï âú òíïãìé äàñïâàóû ðïëïê\nîàšåíô íéñô."
$ws.Cells.Item(21, 5).Value = " Óï, œóï ÿ íïãô óôó ìåçàóû îà\nëàíîÿö... Üóï óïçå âàšà èàòìôãà. Öï-öï-öï!"

$ws.Cells.Item(20, 1).Value = "SCRIPT/P02P01A/us3101.ssb"

# Numeric line-number column
$ws.Cells.Item(19, 2).Value = 18
$ws.Cells.Item(20, 2).Value = 21
$ws.Cells.Item(21, 2).Value = 24

# Row heights
$ws.Rows.Item(19).RowHeight = 57.6
$ws.Rows.Item(20).RowHeight = 43.2
$ws.Rows.Item(21).RowHeight = 21.6

# --- Update the view: scroll position and active cell selection ---
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("C19").Select()
